$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to make the edits, then restore protection
$ws.Unprotect()

# Update the confidential banner text (date 2021-04-26 -> 2021-04-27)
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."
$ws.Range("A41").Value = $newText
$ws.Rows.Item(41).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-38
$ws.Range("D2").Value = 0.0301354361237975
$ws.Range("E2").Value = 0.0004488733279466306
$ws.Range("D3").Value = 0.02989658806710471
$ws.Range("E3").Value = 0.003141361256544517
$ws.Range("D4").Value = 0.03135924255020157
$ws.Range("E4").Value = -0.002341646177262735
$ws.Range("D5").Value = 0.06587645835483159
$ws.Range("E5").Value = 0.002472865943091707
$ws.Range("D6").Value = 0.01513322735900139
$ws.Range("E6").Value = -0.008197977321483529
$ws.Range("D7").Value = 0.01611180871749029
$ws.Range("E7").Value = -0.006836499712147392
$ws.Range("D8").Value = 0.02962643468259296
$ws.Range("E8").Value = -0.001526299311208623
$ws.Range("D9").Value = 0.03365902468830909
$ws.Range("E9").Value = 0.01722356183258689
$ws.Range("D10").Value = 0.02958353478891514
$ws.Range("E10").Value = -0.009308250048990829
$ws.Range("D11").Value = 0.03147731387919776
$ws.Range("E11").Value = 0.01383755908895568
$ws.Range("D12").Value = 0.013868066981034
$ws.Range("E12").Value = -0.01226224482686544
$ws.Range("D13").Value = 0.01469302034189283
$ws.Range("E13").Value = -0.02577794144724732
$ws.Range("D14").Value = 0.01638659992834558
$ws.Range("E14").Value = -0.003632161135875944
$ws.Range("D15").Value = 0.008208179657025012
$ws.Range("E15").Value = 0.0147848196628686
$ws.Range("D16").Value = 0.007212399692512114
$ws.Range("E16").Value = -0.01784422474077652
$ws.Range("D17").Value = 0.03185065889985349
$ws.Range("E17").Value = -0.005308757326085156
$ws.Range("D18").Value = 0.02988692592888898
$ws.Range("E18").Value = -0.001849217638691281
$ws.Range("D19").Value = 0.03181259007528352
$ws.Range("E19").Value = 0.01703872437357634
$ws.Range("D20").Value = 0.02928014364894128
$ws.Range("E20").Value = 0.001748944033790778
$ws.Range("D21").Value = 0.04496256694412463
$ws.Range("E21").Value = -0.008432398978828748
$ws.Range("D22").Value = 0.03319137719866786
$ws.Range("E22").Value = 0.009053330228225631
$ws.Range("D23").Value = 0.03091787607650713
$ws.Range("E23").Value = 0.001875058595580903
$ws.Range("D24").Value = 0.02917656552726868
$ws.Range("E24").Value = 0.01120648545541236
$ws.Range("D25").Value = 0.01512047333655662
$ws.Range("E25").Value = 0.007936507936507908
$ws.Range("D26").Value = 0.01449320732359158
$ws.Range("E26").Value = 0.02200000000000002
$ws.Range("D27").Value = 0.03009311595841261
$ws.Range("E27").Value = 0.003108003108003077
$ws.Range("D28").Value = 0.02962334279836393
$ws.Range("E28").Value = -0.01195073583133277
$ws.Range("D29").Value = 0.03032558700388302
$ws.Range("E29").Value = 0.001605811508315735
$ws.Range("D30").Value = 0.02797382256217489
$ws.Range("E30").Value = 0.003875379939209722
$ws.Range("D31").Value = 0.03589213807272808
$ws.Range("E31").Value = -0.006218503682646426
$ws.Range("D32").Value = 0.03047361096134797
$ws.Range("E32").Value = -0.003278459821428603
$ws.Range("D33").Value = 0.03052771893535604
$ws.Range("E33").Value = -0.002057274522712316
$ws.Range("D34").Value = 0.0305988322726238
$ws.Range("E34").Value = -0.003435558025564611
$ws.Range("D35").Value = 0.02999746079007691
$ws.Range("E35").Value = -0.0001159554730981727
$ws.Range("D36").Value = 0.02909463059519931
$ws.Range("E36").Value = 0.004582890541976692
$ws.Range("D37").Value = 0.03148001927789817
$ws.Range("E37").Value = 0.003584933457742023
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.0008040831423128392

# Restore sheet protection
$ws.Protect()
